$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.443.82"
$ws.Range("E2").Value = '  -0.91%  '

$ws.Range("D3").Value = "'3.757.24"
$ws.Range("E3").Value = '  -0.61%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = "'615.84"
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").Value = "'176.66"
$ws.Range("E6").Value = '  -0.81%  '

$ws.Range("D7").Value = "'3.755.35"
$ws.Range("E7").Value = '  -0.65%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  -1.22%  '

$ws.Range("E10").Value = '  -2.76%  '

$ws.Range("E11").Value = '  +3.14%  '

$ws.Range("D12").Value = "'0.484"
$ws.Range("E12").Value = '  -1.78%  '

$ws.Range("D13").Value = "'40.01"
$ws.Range("E13").Value = '  -2.65%  '

$ws.Range("D14").Value = "'0.0000252"
$ws.Range("E14").Value = '  -3.58%  '

$ws.Range("D15").Value = "'4.387.72"
$ws.Range("E15").Value = '  -0.53%  '

$ws.Range("D16").Value = "'3.758.67"
$ws.Range("E16").Value = '  -0.63%  '

$ws.Range("D17").Value = "'69.554.08"
$ws.Range("E17").Value = '  -0.83%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = "'7.52"
$ws.Range("E18").Value = '  -1.86%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = "'0.120"
$ws.Range("E19").Value = '  -3.22%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'513.29"
$ws.Range("E20").Value = '  +0.41%  '

$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = "'16.43"
$ws.Range("E21").Value = '  -2.01%  '

$ws.Range("D22").Value = "'9.36"
$ws.Range("E22").Value = '  -2.89%  '

$ws.Range("D23").Value = "'0.724"
$ws.Range("E23").Value = '  -0.62%  '

$ws.Range("D24").Value = "'2.52"
$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("D25").Value = "'86.26"
$ws.Range("E25").Value = '  -1.70%  '

$ws.Range("E26").Value = '  -2.93%  '

$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = "'0.0000135"
$ws.Range("E27").Value = '  -3.09%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = "'10.51"
$ws.Range("E28").Value = '  -5.30%  '

$ws.Range("E29").Value = '  +0.21%  '

$ws.Range("D30").Value = "'2.50"
$ws.Range("E30").Value = '  +0.19%  '

$ws.Range("D31").Value = "'2.96"
$ws.Range("E31").Value = '  +3.36%  '

$ws.Range("D32").Value = "'7.88"
$ws.Range("E32").Value = '  +1.80%  '

$ws.Range("D33").Value = "'30.65"
$ws.Range("E33").Value = '  -2.45%  '

$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = '  -1.62%  '

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = '  -0.45%  '

$ws.Range("D37").Value = "'6.13"
$ws.Range("E37").Value = '  -0.95%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = "'0.137"
$ws.Range("E38").Value = '  +2.57%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").Value = "'0.341"
$ws.Range("E39").Value = '  +1.68%  '

$ws.Range("D40").Value = "'445.49"
$ws.Range("E40").Value = '  +5.73%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = "'2.07"
$ws.Range("E41").Value = '  -3.48%  '

$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = "'49.71"
$ws.Range("E42").Value = '  -2.41%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = "'2.93"
$ws.Range("E43").Value = '  +6.01%  '

$ws.Range("D44").Value = "'44.31"
$ws.Range("E44").Value = '  -1.46%  '

$ws.Range("D45").Value = "'8.58"
$ws.Range("E45").Value = '  -2.25%  '

$ws.Range("D46").Value = "'2.946.35"
$ws.Range("E46").Value = '  -3.46%  '

$ws.Range("D47").Value = "'0.0358"
$ws.Range("E47").Value = '  -1.64%  '

$ws.Range("D48").Value = "'27.47"
$ws.Range("E48").Value = '  -0.61%  '

$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = "'138.87"
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("D51").Value = "'2.45"
$ws.Range("E51").Value = '  -1.19%  '
